# repull data, push all data, mean calculation
# Update the dSF column (F) values to reflect the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -1
$ws.Range("F3").Value = -2
$ws.Range("F4").Value = -6
$ws.Range("F7").Value = -4
$ws.Range("F8").Value = 1
$ws.Range("F13").Value = 4
